# Remove the "molFormula" column (H) from the Catalogs_template sheet.
# Deleting the entire column shifts I:O left to H:N and updates the
# sheet's used-range dimension / shared strings automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").EntireColumn.Delete()

# Match the author's final selection (cell L3) recorded in the sheetView.
$ws.Range("L3").Select()
